# Apply the "getting closer on postgame hitter" edits to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 block ---
$ws.Range("M10").Value = ""

# --- Row 12 block ---
$ws.Range("M12").Value = ""

# --- Row 17 block ---
$ws.Range("J17").Value = "FB,CB,CH"

# --- Row 19 block ---
$ws.Range("J19").Value = 3
$ws.Range("M19").Value = ""

# --- Row 20 block ---
$ws.Range("J20").Value = 2

# --- Row 21 block ---
$ws.Range("M21").Value = ""

# --- Row 23 block ---
$ws.Range("J23").Value = "Roblez"
$ws.Range("M23").Value = ""

# --- Row 24 block ---
$ws.Range("M24").Value = "Undefined"

# --- Row 25 block ---
$ws.Range("J25").Value = "88-90 MPH"

# --- Row 26 block ---
$ws.Range("J26").Value = "FB,CB,CH"

# --- Row 28 block ---
$ws.Range("M28").Value = ""

# --- Row 30 block ---
$ws.Range("M30").Value = ""

# --- Row 35 block ---
$ws.Range("J35").Value = "SL,FB,CB,CH"

# --- Row 37 block ---
$ws.Range("J37").Value = 6
$ws.Range("M37").Value = ""

# --- Row 38 block ---
$ws.Range("J38").Value = 0

# --- Row 39 block ---
$ws.Range("M39").Value = ""

# --- Row 41 block ---
$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = "Line Drive"

# --- Row 42 block ---
$ws.Range("M42").Value = "Double"

# --- Row 43 block ---
$ws.Range("J43").Value = "83-85 MPH"

# --- Row 44 block ---
$ws.Range("J44").Value = "SL,FB,CB,CH"

Write-Host "Applied postgame hitter report edits"
